$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price strings (preserve as text, not numbers)
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '27.013.16'
$ws.Range("E2").Value = '  -2.00%  '
$ws.Range("D3").Value = '1.794.48'
$ws.Range("E3").Value = '  -2.32%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = '1.008'
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").Value = '308.56'
$ws.Range("E6").Value = '  -1.74%  '
$ws.Range("D7").Value = '0.4175'
$ws.Range("E7").Value = '  -1.54%  '
$ws.Range("D8").Value = '0.3570'
$ws.Range("E8").Value = '  -2.54%  '
$ws.Range("D9").Value = '0.07047'
$ws.Range("E9").Value = '  -2.56%  '
$ws.Range("D10").Value = '0.8443'
$ws.Range("E10").Value = '  -2.58%  '
$ws.Range("D11").Value = '20.05'
$ws.Range("E11").Value = '  -3.36%  '
$ws.Range("D12").Value = '1.861.86'
$ws.Range("E12").Value = '  -2.12%  '
$ws.Range("D13").Value = '5.257'
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("D14").Value = '6.350'
$ws.Range("E14").Value = '  -2.37%  '
$ws.Range("D15").Value = '0.06852'
$ws.Range("E15").Value = '  -1.67%  '
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("D17").Value = '79.78'
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").Value = '0.000008724'
$ws.Range("E18").Value = '  -3.15%  '
$ws.Range("D19").Value = '1.009'
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").Value = '15.05'
$ws.Range("E20").Value = '  -2.42%  '
$ws.Range("D21").Value = '27.284.24'
$ws.Range("E21").Value = '  -1.43%  '
$ws.Range("D22").Value = '5.046'
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("D23").Value = '10.67'
$ws.Range("E23").Value = '  -1.09%  '
$ws.Range("D24").Value = '2.077.99'
$ws.Range("E24").Value = '  -2.89%  '
$ws.Range("D25").Value = '1.968'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '152.87'
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("D27").Value = '18.18'
$ws.Range("E27").Value = '  -1.26%  '
$ws.Range("D28").Value = '4.990'
$ws.Range("E28").Value = '  -4.91%  '
$ws.Range("D29").Value = '112.74'
$ws.Range("E29").Value = '  -1.71%  '
$ws.Range("D30").Value = '1.664'
$ws.Range("E30").Value = '  -8.72%  '
$ws.Range("D31").Value = '0.08892'
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("D32").Value = '0.7256'
$ws.Range("D33").Value = '2.880'
$ws.Range("E33").Value = '  -2.65%  '
$ws.Range("D34").Value = '4.357'
$ws.Range("E34").Value = '  -4.02%  '
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("D36").Value = '1.081'
$ws.Range("E36").Value = '  -6.43%  '
$ws.Range("E37").Value = '  -2.35%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.05115'
$ws.Range("E38").Value = '  -4.85%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01893'
$ws.Range("E39").Value = '  -2.62%  '
$ws.Range("D40").Value = '0.4960'
$ws.Range("E40").Value = '  -3.05%  '
$ws.Range("D41").Value = '0.1619'
$ws.Range("E41").Value = '  -2.25%  '
$ws.Range("D42").Value = '2.659'
$ws.Range("E42").Value = '  -5.89%  '
$ws.Range("D43").Value = '6.165'
$ws.Range("E43").Value = '  -9.36%  '
$ws.Range("D44").Value = '8.030'
$ws.Range("E44").Value = '  -5.42%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '1.008'
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '10.19'
$ws.Range("E46").Value = '  -2.40%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '104.12'
$ws.Range("E47").Value = '  -0.98%  '
$ws.Range("E48").Value = '  -3.33%  '
$ws.Range("D49").Value = '0.4536'
$ws.Range("E49").Value = '  -3.35%  '
$ws.Range("E50").Value = '  -1.99%  '
$ws.Range("D51").Value = '62.13'
$ws.Range("E51").Value = '  -3.29%  '
